$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset_params")
$ws.Select()
$ws.Rows.Item(13).Resize(2).Insert()
$ws.Range("A13").Value = "cache_rate"
$ws.Range("B13").Value = 0.3
$ws.Range("A14").Value = "ds_type"
$ws.Range("B14").Value = "None"
